# Update column F ("dSF") values for specific rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    7  = 4
    10 = 4
    13 = -1
    23 = 1
    29 = 3
    31 = -5
    39 = 1
    40 = 1
    43 = 2
    44 = -2
    55 = 8
    59 = -2
    62 = 3
    64 = 2
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
